$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '67.269.83'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +7.27%  '
$ws.Range("D3").Value = "'" + '3.584.13'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +3.35%  '
$ws.Range("E4").Value = '  +0.37%  '
$ws.Range("D5").Value = "'" + '416.42'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.57%  '
$ws.Range("D6").Value = "'" + '129.38'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.32%  '
$ws.Range("D7").Value = "'" + '0.651'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.67%  '
$ws.Range("D8").Value = "'" + '3.573.91'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.17%  '
$ws.Range("D9").Value = "'" + '1.00'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.02%  '
$ws.Range("D10").Value = "'" + '0.772'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +6.31%  '
$ws.Range("D11").Value = "'" + '0.176'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +14.33%  '
$ws.Range("D12").Value = "'" + '0.0000332'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +48.37%  '
$ws.Range("D13").Value = "'" + '42.41'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.29%  '
$ws.Range("D14").Value = "'" + '9.87'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.77%  '
$ws.Range("D15").Value = "'" + '4.161.45'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.64%  '
$ws.Range("E16").Value = '  -0.26%  '
$ws.Range("D17").Value = "'" + '20.42'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.85%  '
$ws.Range("D18").Value = "'" + '3.590.84'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.71%  '
$ws.Range("D19").Value = "'" + '1.15'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +6.17%  '
$ws.Range("D20").Value = "'" + '67.159.33'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +7.16%  '
$ws.Range("D21").Value = "'" + '12.28'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.02%  '
$ws.Range("D22").Value = "'" + '451.91'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.67%  '
$ws.Range("D23").Value = "'" + '89.18'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.72%  '
$ws.Range("D24").Value = "'" + '3.15'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.84%  '
$ws.Range("D25").Value = "'" + '13.16'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.01%  '
$ws.Range("D26").Value = "'" + '3.35'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.83%  '
$ws.Range("D27").Value = "'" + '10.10'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.96%  '
$ws.Range("D28").Value = "'" + '34.97'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.81%  '
$ws.Range("E29").Value = '  +1.89%  '
$ws.Range("B30").Value = 'Cosmos'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D30").Value = "'" + '12.38'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.56%  '
$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D31").Value = "'" + '2.73'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.09%  '
$ws.Range("D32").Value = "'" + '0.117'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.46%  '
$ws.Range("D33").Value = "'" + '7.40'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.04%  '
$ws.Range("E34").Value = '  -2.81%  '
$ws.Range("D35").Value = "'" + '41.22'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.98%  '
$ws.Range("E36").Value = '  -0.09%  '
$ws.Range("D37").Value = "'" + '56.69'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.89%  '
$ws.Range("D38").Value = "'" + '0.0494'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.48%  '
$ws.Range("D39").Value = "'" + '0.0₃0743'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +30.32%  '
$ws.Range("E40").Value = '  +9.40%  '
$ws.Range("B41").Value = 'FirstDigitalUSD'
$ws.Range("C41").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D41").Value = "'" + '0.998'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.01%  '
$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").Value = "'" + '3.07'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.56%  '
$ws.Range("E43").Value = '  +1.47%  '
$ws.Range("D44").Value = "'" + '149.01'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.01%  '
$ws.Range("D45").Value = "'" + '0.315'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.23%  '
$ws.Range("D46").Value = "'" + '3.27'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.58%  '
$ws.Range("D47").Value = "'" + '4.33'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.79%  '
$ws.Range("D48").Value = "'" + '1.97'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.49%  '
$ws.Range("E49").Value = '  -3.16%  '
$ws.Range("B50").Value = 'BitcoinSV'
$ws.Range("C50").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D50").Value = "'" + '115.58'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +6.11%  '
$ws.Range("B51").Value = 'Celestia'
$ws.Range("C51").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D51").Value = "'" + '15.64'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.68%  '

